$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific Price (column D) cells that would otherwise be auto-detected
# as numbers to stay as text, matching the original inline-string storage.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values from the crypto data refresh.
$ws.Range("D2").Value = "61.283.34"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "3.374.71"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "571.26"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "138.15"
$ws.Range("E6").Value = "  +8.89%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.376.04"
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").Value = "7.60"
$ws.Range("E10").Value = "  +5.25%  "
$ws.Range("E11").Value = "  +4.59%  "
$ws.Range("E12").Value = "  +4.64%  "
$ws.Range("D13").Value = "3.941.86"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("E14").Value = "  +2.35%  "
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("D16").Value = "3.369.87"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "61.337.29"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("E19").Value = "  +6.01%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "9.45"
$ws.Range("E20").Value = "  +4.02%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "5.80"
$ws.Range("E21").Value = "  +3.47%  "
$ws.Range("D22").Value = "382.08"
$ws.Range("E22").Value = "  +8.51%  "
$ws.Range("E23").Value = "  +3.78%  "
$ws.Range("D24").Value = "3.507.38"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "70.84"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("E27").Value = "  +9.81%  "
$ws.Range("D28").Value = "1.66"
$ws.Range("E28").Value = "  +13.04%  "
$ws.Range("E29").Value = "  +8.99%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  +3.01%  "
$ws.Range("E32").Value = "  +4.80%  "
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "3.399.98"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("D36").Value = "23.48"
$ws.Range("E36").Value = "  +4.99%  "
$ws.Range("D37").Value = "5.57"
$ws.Range("E37").Value = "  +2.37%  "
$ws.Range("D38").Value = "6.96"
$ws.Range("E38").Value = "  +3.37%  "
$ws.Range("E39").Value = "  +4.19%  "
$ws.Range("D40").Value = "163.80"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").Value = "0.0804"
$ws.Range("E41").Value = "  +6.54%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "41.56"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "4.41"
$ws.Range("E44").Value = "  +4.07%  "
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("E46").Value = "  +7.55%  "
$ws.Range("E47").Value = "  +5.77%  "
$ws.Range("D48").Value = "23.29"
$ws.Range("E48").Value = "  +3.66%  "
$ws.Range("E49").Value = "  +5.21%  "
$ws.Range("D50").Value = "23.20"
$ws.Range("E50").Value = "  +11.80%  "
$ws.Range("D51").Value = "2.43"
$ws.Range("E51").Value = "  +10.98%  "
